$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = "Lab 0: Getting started w/ Jupyter notebook / test submitting a lab"
$ws.Range("G2").Value = "Hwk 1: Collecting personal network survey data"

# Row 3
$ws.Range("F3").Value = "Lab 1: Analyzing personal network data; review of bootstrap"
$ws.Range("G3").Value = ""

# Row 5
$ws.Range("F5").Value = "Lab 2: Getting started with complete network data"
$ws.Range("G5").Value = "Hwk 2: Analyzing personal network data"

# Row 7
$ws.Range("G7").Value = "Hwk 3: Complete network data"

# Row 9
$ws.Range("F9").Value = "Lab 4 - Exploring mathematical models"
$ws.Range("G9").Value = "Hwk 4: Problem set I"

# Row 11
$ws.Range("F11").Value = "Lab 5 - Two-mode networks; the friendship paradox"
$ws.Range("G11").Value = "Hwk 5: Advanced complete network data"

# Row 13
$ws.Range("F13").Value = "Midterm review / question session"
$ws.Range("G13").Value = "Hwk 6: Problem set II"

# Row 14
$ws.Range("D14").Value = "Midterm review"
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""

# Row 15
$ws.Range("D15").Value = "Midterm"

# Row 16
$ws.Range("D16").Value = "Community detection"
$ws.Range("E16").Value = "Community detection"

# Row 17
$ws.Range("D17").Value = "Empirical studies of network structure"

# Row 18
$ws.Range("F18").Value = ""

# Row 19
$ws.Range("F19").Value = "Lab 6: Simple contagion"
$ws.Range("G19").Value = "Hwk 7: Centrality and the SIR model"

# Row 20
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""

# Row 21
$ws.Range("F21").Value = "Lab 7: Mini-project prep"
$ws.Range("G21").Value = "Mini-project"

# Row 22
$ws.Range("F22").Value = ""

# Row 23
$ws.Range("F23").Value = "Mini-project support"

# Row 24
$ws.Range("G24").Value = ""

# Row 25
$ws.Range("G25").Value = "Hwk 8: Problem set III"

# Row 29
$ws.Range("F29").Value = "Problem set + mini-project support"
